$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Title 1")

# Replace the whole multi-run title ("<space>Notebook Tracking <br/>System")
# with a single run containing the new title text, keeping the existing
# Sitka Banner / Segoe UI + size/bold formatting of the first run.
$tr = $sh.TextFrame.TextRange
$tr.Delete()
$sh.TextFrame.TextRange.Text = "Enfocar"
